$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")
$ws.Activate()

# --- Mark a batch of previously "Not Passed" challenges as "Passed" ---
# (Row 86-93 edited first, then row 80 - matches the order the review was done in.)
$ws.Range("D86").Value = "PassedUpdating Object Properties"
$ws.Range("D87").Value = "PassedAdd New Properties to a JavaScript Object"
$ws.Range("D88").Value = "PassedDelete Properties from a JavaScript Object"
$ws.Range("D89").Value = "PassedUsing Objects for Lookups"
$ws.Range("D90").Value = "PassedTesting Objects for Properties"
$ws.Range("D91").Value = "PassedManipulating Complex Objects"
$ws.Range("D92").Value = "PassedAccessing Nested Objects"
$ws.Range("D93").Value = "PassedAccessing Nested Arrays"
$ws.Range("D80").Value = "PassedReturn Early Pattern for Functions"

# --- Highlight the next block of challenges (still pending review) in yellow ---
$ws.Range("D94:D136").Interior.Color = 65535

# --- Scroll the frozen pane down to the section currently under review ---
$win = $excel.ActiveWindow
$win.ScrollRow = 122
$win.ScrollColumn = 1
$ws.Range("D94:D136").Select()
